$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-recognized by Excel
# as a number must be forced to Text format first, so they stay inline/shared
# strings (matching the source data which stores these as text).

# Row 2
$ws.Range("D2").Value = "63.880.98"
$ws.Range("E2").Value = "  +2.65%  "

# Row 3
$ws.Range("D3").Value = "3.053.85"
$ws.Range("E3").Value = "  +1.90%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.67"
$ws.Range("E5").Value = "  +2.41%  "

# Row 6
$ws.Range("E6").Value = "  +2.78%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "3.052.84"
$ws.Range("E8").Value = "  +1.99%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  +4.65%  "

# Row 10
$ws.Range("E10").Value = "  +5.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.05"
$ws.Range("E11").Value = "  -12.22%  "

# Row 12
$ws.Range("E12").Value = "  +7.51%  "

# Row 13
$ws.Range("E13").Value = "  +4.94%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.26"
$ws.Range("E14").Value = "  +3.89%  "

# Row 15
$ws.Range("D15").Value = "3.551.49"
$ws.Range("E15").Value = "  +2.30%  "

# Row 16
$ws.Range("D16").Value = "63.943.40"
$ws.Range("E16").Value = "  +2.63%  "

# Row 17
$ws.Range("D17").Value = "3.056.89"
$ws.Range("E17").Value = "  +1.99%  "

# Row 18
$ws.Range("E18").Value = "  +1.82%  "

# Row 19
$ws.Range("E19").Value = "  +2.89%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.57"
$ws.Range("E20").Value = "  +1.83%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.95"
$ws.Range("E21").Value = "  +4.27%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("E22").Value = "  +4.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.62"
$ws.Range("E23").Value = "  +6.25%  "

# Row 24
$ws.Range("E24").Value = "  +13.99%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.05"
$ws.Range("E25").Value = "  +3.48%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("E27").Value = "  +2.52%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.93"
$ws.Range("E28").Value = "  +3.91%  "

# Row 29
$ws.Range("E29").Value = "  +1.57%  "

# Row 30
$ws.Range("E30").Value = "  +0.00%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.29"
$ws.Range("E31").Value = "  +3.64%  "

# Row 32
$ws.Range("E32").Value = "  +1.42%  "

# Row 33
$ws.Range("E33").Value = "  +3.66%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.66"
$ws.Range("E34").Value = "  +2.04%  "

# Row 35
$ws.Range("E35").Value = "  +6.57%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.76"
$ws.Range("E36").Value = "  +0.58%  "

# Row 37
$ws.Range("E37").Value = "  +3.76%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "443.75"
$ws.Range("E38").Value = "  -1.84%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0807"
$ws.Range("E39").Value = "  -0.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  +12.07%  "

# Row 41
$ws.Range("D41").Value = "2.990.85"
$ws.Range("E41").Value = "  +1.66%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.24"
$ws.Range("E42").Value = "  +2.40%  "

# Row 43
$ws.Range("E43").Value = "  +0.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.79"
$ws.Range("E44").Value = "  +3.84%  "

# Row 45
$ws.Range("E45").Value = "  +5.10%  "

# Row 46
$ws.Range("E46").Value = "  +7.89%  "

# Row 48
$ws.Range("E48").Value = "  +4.00%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.25"
$ws.Range("E49").Value = "  +2.69%  "

# Row 50
$ws.Range("D50").Value = "0.0₃0517"
$ws.Range("E50").Value = "  +4.95%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.07"
$ws.Range("E51").Value = "  +2.96%  "
